$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting templates before underlying cells change ---
$ws.Range("L4").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("L9").PasteSpecial(-4122)

$ws.Range("J9").Copy()
$ws.Range("J10").PasteSpecial(-4122)

$ws.Range("J2").Copy()
$ws.Range("J9").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Clear cells that are removed entirely ---
$ws.Range("L4").Clear()
$ws.Range("O5").Clear()
$ws.Range("P5").Clear()
$ws.Range("Q5").Clear()
$ws.Range("R5").Clear()
$ws.Range("S5").Clear()
$ws.Range("Q6").Clear()
$ws.Range("M8").Clear()
$ws.Range("N8").Clear()

# --- Set new / changed values ---
$ws.Range("A1").Value = 's'
$ws.Range("J2").Value = 'אוראל'
$ws.Range("M2").Value = 'f'
$ws.Range("N2").Value = 'b.c'
$ws.Range("O2").Value = 'f'
$ws.Range("P2").Value = 'a.b.c'
$ws.Range("U2").Value = 4
$ws.Range("A3").Value = 'עמדה 1'
$ws.Range("J3").Value = 'זיו שהינו'
$ws.Range("M3").Value = 'f'
$ws.Range("N3").Value = 'b.c'
$ws.Range("O3").Value = 'f'
$ws.Range("P3").Value = 'a.b.c'
$ws.Range("Q3").Value = 'a'
$ws.Range("S3").Value = 'a.b'
$ws.Range("U3").Value = 4
$ws.Range("J4").Value = 'רוני'
$ws.Range("M4").Value = 'f'
$ws.Range("N4").Value = 'b.c'
$ws.Range("O4").Value = 'f'
$ws.Range("P4").Value = 'a.b.c'
$ws.Range("Q4").Value = 'a'
$ws.Range("U4").Value = 3
$ws.Range("J5").Value = 'יניב שטיינר'
$ws.Range("L5").Value = 'חתונה'
$ws.Range("M5").Value = 'a'
$ws.Range("N5").Value = 'a'
$ws.Range("U5").Value = 0
$ws.Range("J6").Value = 'אור'
$ws.Range("L6").Value = 'חתונה ברביעי'
$ws.Range("M6").Value = 'a.b.c'
$ws.Range("P6").Value = 'a'
$ws.Range("R6").Value = 'b.c'
$ws.Range("S6").Value = 'a.c'
$ws.Range("U6").Value = 4
$ws.Range("A7").Value = 'עמדה 2'
$ws.Range("J7").Value = 'שבת'
$ws.Range("L7").Value = 'עד 4 משמרות לבינתיים 🫡'
$ws.Range("M7").Value = 'a.b.c'
$ws.Range("O7").Value = 'a'
$ws.Range("P7").Value = 'b.c'
$ws.Range("Q7").Value = 'a.b.c'
$ws.Range("R7").Value = 'a.b.c'
$ws.Range("S7").Value = 'a.b.c'
$ws.Range("U7").Value = 5
$ws.Range("J8").Value = 'עמית בלסן'
$ws.Range("O8").Value = 'c'
$ws.Range("R8").Value = 'a'
$ws.Range("S8").Value = 'c'
$ws.Range("U8").Value = 2
$ws.Range("J9").Value = 'אמור'
$ws.Range("L9").Value = 'חינה בשלישי בערב. 
אם אפשר בלי : (b)-(b) בשישי שבת. 
ומזל טוב לשטיינר'
$ws.Range("M9").Value = 'a.b'
$ws.Range("N9").Value = 'b.c'
$ws.Range("P9").Value = 'a.b.c'
$ws.Range("Q9").Value = 'a'
$ws.Range("R9").Value = 'a.b.c'
$ws.Range("S9").Value = 'a.b'
$ws.Range("T9").Value = 6
$ws.Range("U9").Value = 5
$ws.Range("J10").Value = 'תגבור'

# --- Restore row 9 height (Excel auto-grows it for the multi-line L9 comment) ---
$ws.Rows.Item(9).EntireRow.AutoFit()
